$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D28").Value = "WSL Docker 설치"
$ws.Range("E28").Value = "https://ropiens.tistory.com/160"

$ws.Range("D46").Value = "[Bioinformatics] 2021년 10월, 유전체 빅데이터 통합 분석 플랫폼 Bio-Express 활용 교육"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/416"
